$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the header text in A1 ("name") and D1 ("geneid")
$a1 = $ws.Range("A1").Value2
$d1 = $ws.Range("D1").Value2
$ws.Range("A1").Value = $d1
$ws.Range("D1").Value = $a1

# Re-apply the "Normal" style on the header row, except C1 which keeps its
# original (untouched) style
$ws.Range("A1").Style = "Normal"
$ws.Range("B1").Style = "Normal"
$ws.Range("D1").Style = "Normal"
$ws.Range("E1").Style = "Normal"

# Row 1 height
$ws.Rows(1).RowHeight = 12.8

# Final selection
$ws.Range("C1").Select()
